$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# ---------------------------------------------------------------------------
# 1. Version string in the title cell (A1 / shared string 0)
# ---------------------------------------------------------------------------
$ws.Range("A1").Value = $ws.Range("A1").Text.Replace( `
    "Kwaliteitsaanpak ICTU Software Realisatie versie 1.2.1-build.4, 29-08-2018.", `
    "Kwaliteitsaanpak ICTU Software Realisatie versie 1.2.1-build.37, 31-08-2018.")

# ---------------------------------------------------------------------------
# 2. Product names in the table (B5 / B6)
# ---------------------------------------------------------------------------
$ws.Range("B5").Value = "1. BIA (Business impact analysis)*"
$ws.Range("B6").Value = "2. PIA (Privacy impact analysis)**"

# ---------------------------------------------------------------------------
# 3. Comment on B4 - M01 Op te leveren producten
# ---------------------------------------------------------------------------
$c = $ws.Range("B4").Comment
$t = $c.Text()
$t = $t.Replace( `
    "| Product                 | Voorbereidings᠆fase  | Voorbereidings᠆fase met onderzoek  | Realisatie᠆fase |", `
    "| Product                | Voorbereidings᠆fase  | Voorbereidings᠆fase met onderzoek  | Realisatie᠆fase |")
$t = $t.Replace( `
    "| Business impact analysis (BIA)*                                                        | ✔ | ✔ | ✔ |", `
    "| BIA (Business impact analysis)*                                                        | ✔ | ✔ | ✔ |")
$t = $t.Replace( `
    "| Privacy impact analysis (PIA)**                                                        | ✔ | ✔ | ✔ |", `
    "| PIA (Privacy impact analysis)**                                                        | ✔ | ✔ | ✔ |")
$t = $t.Replace( `
    "- De beschrijving van niet-functionele eisen is gebaseerd op ISO-25010, BIR en SSD, en bevat een prioritering van de niet-functionele eisen. De beschrijving van niet-functionele eisen is gebaseerd op het ICTU NFE-template. De beschrijving bevat in ieder geval eisen aan toegangsbeveiliging, aan beheerfuncties, aan logging en aan het gewenste gedrag van de software bij uitval van infrastructurele diensten, zoals een log-server;", `
    "- De beschrijving van niet-functionele eisen is gebaseerd op ISO (International Organization for Standardization)-25010, de BIR (Baseline Informatiebeveiliging Rijksdienst) en de methode Grip op SSD (Secure software development) van het CIP (Centrum Informatiebeveiliging en Privacybescherming), en bevat een prioritering van de niet-functionele eisen. De beschrijving van niet-functionele eisen is gebaseerd op het ICTU NFE (Niet-functionele eisen)-template. De beschrijving bevat in ieder geval eisen aan toegangsbeveiliging, aan beheerfuncties, aan logging en aan het gewenste gedrag van de software bij uitval van infrastructurele diensten, zoals een log-server;")
$t = $t.Replace( `
    "- De ontwerp- en architectuurdocumentatie bestaat uit een projectstartarchitectuur (PSA), een softwarearchitectuurdocument (SAD), een infrastructuurarchitectuur (IA), een globaal functioneel ontwerp (GFO) bijvoorbeeld in de vorm van use cases, en een prototype en/of interactieontwerp. De architectuurdocumenten moeten expliciet inzichtelijk maken hoe aan de niet-functionele eisen wordt voldaan door uit te werken welke oplossingen en mechanieken gekozen zijn, bijvoorbeeld voor identificatie, authenticatie, autorisatie, concurrency, transactionele verwerking of logging;", `
    "- De ontwerp- en architectuurdocumentatie bestaat uit een PSA (Projectstartarchitectuur), een SAD (Softwarearchitectuurdocument), een IA (Infrastructuurarchitectuur), een GFO (Globaal functioneel ontwerp) bijvoorbeeld in de vorm van use cases, en een prototype en/of interactieontwerp. De architectuurdocumenten moeten expliciet inzichtelijk maken hoe aan de niet-functionele eisen wordt voldaan door uit te werken welke oplossingen en mechanieken gekozen zijn, bijvoorbeeld voor identificatie, authenticatie, autorisatie, concurrency, transactionele verwerking of logging;")
$t = $t.Replace( `
    "- De testdocumentatie bestaat uit een mastertestplan, gemaakt op basis van een productrisicoanalyse (PRA). Beveiligingstesten zijn een integraal onderdeel van het mastertestplan en worden als zodanig afgestemd met de opdrachtgever;", `
    "- De testdocumentatie bestaat uit een mastertestplan, gemaakt op basis van een PRA (Productrisicoanalyse). Beveiligingstesten zijn een integraal onderdeel van het mastertestplan en worden als zodanig afgestemd met de opdrachtgever;")
$t = $t.Replace( `
    "- Het informatiebeveiligingsplan is gebaseerd op een dreigingen- en kwetsbaarhedenanalyse (TVA, threat and vulnerability assessment) en bevat een maatregelenselectie informatiebeveiliging. De TVA wordt tijdens de voorfase opgesteld op basis van de resultaten van de BIA, de eventuele PIA en inhoud van de ontwerp- en architectuurdocumentatie. Een TVA levert een deel van een traceerbare onderbouwing voor de te treffen beveiligingsmaatregelen;", `
    "- Het informatiebeveiligingsplan is gebaseerd op een dreigingen- en kwetsbaarhedenanalyse (TVA (Threat and vulnerability assessment)) en bevat een maatregelenselectie informatiebeveiliging. De TVA wordt tijdens de voorfase opgesteld op basis van de resultaten van de BIA, de eventuele PIA en inhoud van de ontwerp- en architectuurdocumentatie. Een TVA levert een deel van een traceerbare onderbouwing voor de te treffen beveiligingsmaatregelen;")
$c.Text($t)

# ---------------------------------------------------------------------------
# 4. Comment on B24 - M13 Gebruik van ISO-25010
# ---------------------------------------------------------------------------
$c = $ws.Range("B24").Comment
$t = $c.Text()
$t = $t.Replace( `
    "De standaard ISO/IEC 25010:2011, kortweg ""ISO-25010"", biedt een model voor het beschrijven van productkwaliteit.", `
    "De standaard ISO/IEC (International Electrotechnical Commission) 25010:2011, kortweg ""ISO-25010"", biedt een model voor het beschrijven van productkwaliteit.")
$c.Text($t)

# ---------------------------------------------------------------------------
# 5. Comment on B25 - M02 Continu voldoen aan kwaliteitsnormen
# ---------------------------------------------------------------------------
$c = $ws.Range("B25").Comment
$t = $c.Text()
$t = $t.Replace( `
    "Tijdens de realisatiefase van softwarerealisatieprojecten wordt het voldoen aan de kwaliteitsnormen diverse malen per uur gemeten door het 'Kwaliteitssysteem' (HQ).", `
    "Tijdens de realisatiefase van softwarerealisatieprojecten wordt het voldoen aan de kwaliteitsnormen diverse malen per uur gemeten door het 'Kwaliteitssysteem', genaamd HQ (Holistic Software Quality Reporting).")
$c.Text($t)

# ---------------------------------------------------------------------------
# 6. Comment on B28 - M26 Periodieke beoordeling informatiebeveiliging
# ---------------------------------------------------------------------------
$c = $ws.Range("B28").Comment
$t = $c.Text()
$t = $t.Replace( `
    "Overheidsspecifieke beveiligingsnormen of -raamwerken, zoals de Baseline Informatiebeveiliging Rijksdienst (BIR), bieden een basis voor de beoordeling.", `
    "Overheidsspecifieke beveiligingsnormen of -raamwerken, zoals de BIR (Baseline Informatiebeveiliging Rijksdienst), bieden een basis voor de beoordeling.")
$t = $t.Replace( `
    "of bekende kwetsbaarheden (OWASP) vermeden zijn", `
    "of bekende kwetsbaarheden (zoals bijvoorbeeld in de OWASP Top 10 genoemd) vermeden zijn")
$c.Text($t)

# ---------------------------------------------------------------------------
# 7. Comment on B36 - M07 Continuous delivery pipeline
# ---------------------------------------------------------------------------
$c = $ws.Range("B36").Comment
$t = $c.Text()
$t = $t.Replace( `
    "ICTU gebruikt Jenkins of Team Foundation Server (TFS) als tool voor de implementatie van de continuous delivery pipeline.", `
    "ICTU gebruikt Jenkins of TFS (Team Foundation Server) als tool voor de implementatie van de continuous delivery pipeline.")
$c.Text($t)

# ---------------------------------------------------------------------------
# 8. Comment on B45 - M08 Technische schuld
# ---------------------------------------------------------------------------
$c = $ws.Range("B45").Comment
$t = $c.Text()
$t = $t.Replace( `
    "ICTU gebruikt [HQ](https://github.com/ICTU/quality-report/) (een door ICTU ontwikkeld, open source, geautomatiseerd kwaliteitssysteem) om bestaande technische schuld inzichtelijk te maken", `
    "ICTU gebruikt HQ (Holistic Software Quality Reporting), een door ICTU ontwikkeld, open source, geautomatiseerd kwaliteitssysteem, om bestaande technische schuld inzichtelijk te maken")
$c.Text($t)

# ---------------------------------------------------------------------------
# 9. Comment on B55 - M15 Open source tools
# ---------------------------------------------------------------------------
$c = $ws.Range("B55").Comment
$t = $c.Text()
$t = $t.Replace( `
    "Conform de rationale uit NORA voor het gebruik van open source tools,", `
    "Conform de rationale uit NORA (Nederlandse Overheid Referentiearchitectuur) voor het gebruik van open source tools,")
$c.Text($t)

# ---------------------------------------------------------------------------
# 10. Comment on B56 - M16 Verplichte tools
# ---------------------------------------------------------------------------
$c = $ws.Range("B56").Comment
$t = $c.Text()
$t = $t.Replace( `
    "2. Jenkins voor Javaprojecten en Team Foundation Server (TFS) voor DotNet-projecten,", `
    "2. Jenkins voor Javaprojecten en TFS (Team Foundation Server) voor DotNet-projecten,")
$t = $t.Replace( `
    "5. Reporting (Birt),", `
    "5. Reporting (ontwikkeld met behulp van BIRT (Business Intelligence Reporting Tool)),")
$t = $t.Replace( `
    "7. OpenVAS en OWASP ZAP,", `
    "7. OpenVAS (Vulnerability Assessment System) en OWASP (Open Web Application Security Project) ZAP (Zed Attack Proxy),")
$c.Text($t)

# ---------------------------------------------------------------------------
# 11. Comment on B71 - M19 Afgeschermde digitale omgeving
# ---------------------------------------------------------------------------
$c = $ws.Range("B71").Comment
$t = $c.Text()
$t = $t.Replace( `
    "ICTU ondersteunt dit met Docker en/of virtuele machines (VM) en een VLAN per project.", `
    "ICTU ondersteunt dit met Docker en/of virtuele machines en een VLAN (Virtual local area network) per project.")
$c.Text($t)

# ---------------------------------------------------------------------------
# 12. Shared strings used in the "Verplichte tools" (M16) table rows
# ---------------------------------------------------------------------------
$ws.Range("B58").Value = $ws.Range("B58").Text.Replace( `
    "Jenkins voor Javaprojecten en Team Foundation Server (TFS) voor DotNet-projecten,", `
    "Jenkins voor Javaprojecten en TFS (Team Foundation Server) voor DotNet-projecten,")
$ws.Range("B61").Value = $ws.Range("B61").Text.Replace( `
    "Reporting (Birt),", `
    "Reporting (ontwikkeld met behulp van BIRT (Business Intelligence Reporting Tool)),")
$ws.Range("B63").Value = $ws.Range("B63").Text.Replace( `
    "OpenVAS en OWASP ZAP,", `
    "OpenVAS (Vulnerability Assessment System) en OWASP (Open Web Application Security Project) ZAP (Zed Attack Proxy),")
